$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q4: numeric-looking code value must remain text, so set as text explicitly.
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "13005000"

$ws.Range("Q5").Value = "Merici College"

$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "52.0411"

$ws.Range("Q7").Value = "Computer Programming, Specific Applications"
